$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 459, pushing existing rows 459:478 down to 460:479
$ws.Rows.Item(459).EntireRow.Insert()

# Populate the newly inserted row 459 with the new record
$ws.Cells.Item(459, 1).Value = 10
$ws.Cells.Item(459, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(459, 3).Value = "La Araucanía"
$ws.Cells.Item(459, 4).Value = 45008
$ws.Cells.Item(459, 5).Value = 9
$ws.Cells.Item(459, 6).Value = 100112009
$ws.Cells.Item(459, 7).Value = "Acelga"
$ws.Cells.Item(459, 8).Value = "Sin especificar"
$ws.Cells.Item(459, 9).Value = "Primera"
$ws.Cells.Item(459, 10).Value = 65
$ws.Cells.Item(459, 11).Value = 8000
$ws.Cells.Item(459, 12).Value = 8000
$ws.Cells.Item(459, 13).Value = 8000
$ws.Cells.Item(459, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(459, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(459, 16).Value = 667
$ws.Cells.Item(459, 17).Value = 12
$ws.Cells.Item(459, 18).Value = "Hortaliza"
